# Add reply-length [B] values for the Timepix commands on the
# "all_systems" sheet (column K, rows 4-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("all_systems")

$ws.Range("K4").Value = "0x0c"
$ws.Range("K5").Value = "0x03"
$ws.Range("K6").Value = "0x08"
$ws.Range("K7").Value = "0x08"
$ws.Range("K8").Value = "0x08"
$ws.Range("K9").Value = "0x04"

# Leave the selection on K10 (matches the saved view state in the commit).
$ws.Range("K10").Select() | Out-Null
